$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price rows (2, 4, 5) get their Fecha/Volumen/Precio values
# rotated: row2 <- old row4, row4 <- old row5, row5 <- old row2.

$ws.Range("D2").Value = 44714
$ws.Range("J2").Value = 80

$ws.Range("D4").Value = 44804
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 9500
$ws.Range("M4").Value = 9750
$ws.Range("P4").Value = 542

$ws.Range("D5").Value = 44792
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 9000
$ws.Range("M5").Value = 9500
$ws.Range("P5").Value = 528
